$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (+556293410347 / 62 / 2024-09-30) entirely; rows below shift up.
$ws.Rows("3").Delete()
